$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.076.55"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "3.120.76"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.36"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.53"
$ws.Range("E6").Value = "  +3.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.116.08"
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.42"
$ws.Range("E10").Value = "  -3.45%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").Value = "3.636.47"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").Value = "67.056.42"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "3.122.97"
$ws.Range("E19").Value = "  +1.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.19"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "478.14"
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.71"
$ws.Range("E23").Value = "  +3.48%  "
$ws.Range("E24").Value = "  +3.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.94"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.01"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.43"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.99"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").Value = "0.0₃0974"
$ws.Range("E33").Value = "  -4.42%  "
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.978"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.72"
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("E39").Value = "  +2.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.06"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").Value = "2.813.65"
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0357"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "381.25"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.59"
$ws.Range("E47").Value = "  -10.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.88"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("E51").Value = "  -0.45%  "
